$d = $word.ActiveDocument
$sec = $d.Sections.First

# This section has two distinct headers and two distinct footers:
#   - a "first page" header/footer  (WdHeaderFooterIndex = 2)
#   - a "default/primary" header/footer (WdHeaderFooterIndex = 1)
# Each one contains a single inline picture - the BTEC logo in the headers,
# the Pearson logo in the footers - whose wp:docPr (InlineShape.Name) needs
# renaming.

# --- Footers: Pearson logo "image2.png" -> "image1.png" ---
# Selecting the picture's range first (instead of renaming the InlineShape
# straight off the HeaderFooter.Range collection) makes sure the rename
# lands on a freshly-addressed object.
$ftrFirst = $sec.Footers.Item(2)
$ftrFirst.Range.InlineShapes.Item(1).Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

$ftrDefault = $sec.Footers.Item(1)
$ftrDefault.Range.InlineShapes.Item(1).Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# --- Headers: BTEC logo "image1.jpg" -> "image2.jpg" ---
$hdrFirst = $sec.Headers.Item(2)
$hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"

$hdrDefault = $sec.Headers.Item(1)
$hdrDefault.Range.InlineShapes.Item(1).Name = "image2.jpg"
